$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D (Ref-Genome) to fit the new data
$ws.Columns("D").ColumnWidth = 19

# New BOWTIE mapping results (rows 14-15): mapping time + SAM file sizes
$ws.Range("G14").Value = 156.96
$ws.Range("I14").Value = "129M"
$ws.Range("H14").Value = "121M"

$ws.Range("G15").Value = 277.3
$ws.Range("I15").Value = "257MB"
$ws.Range("H15").Value = "241MB"

# Ref. Genome size, now filled in across the whole table
$ws.Range("E14").Value = "128M"
$ws.Range("E15").Value = "128M"
$ws.Range("E16").Value = "128M"

$ws.Range("E10").Value = "163MB"
$ws.Range("E11").Value = "163MB"
$ws.Range("E12").Value = "163MB"

# Row 16 (BOWTIE, 400MB) mapping results
$ws.Range("G16").Value = 497.73
$ws.Range("I16").Value = "504MB"
$ws.Range("H16").Value = "472MB"

$ws.Range("H17").Select()
